$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Identifying the Human Rights Impacts of Palm Oil', 'Dr. Maria Schmidt', 'A, B', 8.6, '2025-06-25 16:28:08')
    ,@('Handlungsleitfaden zur Durchführung der Risikoanalyse für Kakao produzierende Länder', 'Prof. Dr. James Allen', 'Undefined', 6.2, '2025-06-25 16:28:08')
    ,@('UNEP FI Human Rights Guidance Tool for the Financial Sector "Infrastructure"', 'Abdallah Reyati', 'K', 8.2, '2025-06-25 16:28:08')
    ,@('Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries', 'Lisa Müller', 'C13', 5.4, '2025-06-25 16:28:08')
    ,@('Praxislotse Wirtschaft & Menschenrechte', 'Dr. Maria Schmidt', 'Undefined', 9, '2025-06-25 16:28:08')
    ,@('UNEP FI Human Rights Guidance Tool for the Financial Sector "Oil and Gas"', 'Prof. Dr. James Allen', 'B, K', 6, '2025-06-25 16:28:08')
    ,@('UNEP FI Human Rights Guidance Tool for the Financial Sector "Utilities and Waste Management"', 'Abdallah Reyati', 'D, K', 8.9, '2025-06-25 16:28:08')
    ,@('Strengthening Protection Against Trafficking in Persons in Federal and Corporate Supply Chains', 'Lisa Müller', 'Undefined', 5, '2025-06-25 16:28:08')
    ,@('Human Rights Risks in Mining, A Baseline Study', 'Dr. Maria Schmidt', 'B', 8.8, '2025-06-25 16:28:08')
    ,@('2022 List Of Goods Produced By Child Labor Or Forced Labor', 'Prof. Dr. James Allen', 'Undefined', 5.9, '2025-06-25 16:28:08')
    ,@('Die Achtung von Menschenrechten entlang globaler Wertschöpfungskette', 'Abdallah Reyati', 'Undefined', 8.5, '2025-06-25 16:28:08')
    ,@('Global Estimates of Modern Slavery, Forced Labour and Forced Marriage', 'Lisa Müller', 'Undefined', 4.8, '2025-06-25 16:28:08')
    ,@('Low Prices Drive Natural Rubber Producers Into Poverty', 'Dr. Maria Schmidt', 'Undefined', 8.6, '2025-06-25 16:28:08')
    ,@('Global Mercury Assessment 2018', 'Prof. Dr. James Allen', 'Undefined', 6.2, '2025-06-25 16:28:08')
    ,@('Sustainability and Circularity in the Textile Value Chain - A Global Roadmap', 'Abdallah Reyati', 'Undefined', 8.2, '2025-06-25 16:28:08')
    ,@('Forced and Child Labour in the Cotton Industry', 'Lisa Müller', 'Undefined', 5.8, '2025-06-25 16:28:08')
    ,@('Forced, child and trafficked labour in the cocoa industry', 'Dr. Maria Schmidt', 'Undefined', 8.6, '2025-06-25 16:28:08')
    ,@('Ermittlung von potentiell POP-haltigen Abfällen und Recyclingstoffen - Ableitung von Grenzwerten', 'Prof. Dr. James Allen', 'Undefined', 6, '2025-06-25 16:28:08')
    ,@('Kurzzeitige Chlorparafine (SCCP) Vorkommen, Verwendung und Rechtssetzung zu kurzketigen Chlorparaffinen in Produkten und Abfällen', 'Abdallah Reyati', 'Undefined', 8.1, '2025-06-25 16:28:08')
    ,@('Hexabromocyclododecane', 'Lisa Müller', 'Undefined', 5, '2025-06-25 16:28:08')
    ,@('Fairtrade Risk Map', 'Dr. Maria Schmidt', 'Undefined', 8.4, '2025-06-25 16:28:08')
    ,@('Typical Wastes Generated By Industry Sector', 'Prof. Dr. James Allen', 'Undefined', 6.6, '2025-06-25 16:28:08')
    ,@('Cocoa Barometer 2020', 'Abdallah Reyati', 'Undefined', 8.1, '2025-06-25 16:28:08')
    ,@('Handlungsleitfaden zur Durchführung der Risikoanalysen für Kakao produzierende Länder', 'Lisa Müller', 'Undefined', 4.8, '2025-06-25 16:28:08')
    ,@('Human Rights Toolkit for Financial Institutions', 'Dr. Maria Schmidt', 'K', 9.4, '2025-06-25 16:28:08')
    ,@('Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries', 'Prof. Dr. James Allen', 'C13', 5.5, '2025-06-25 16:28:08')
    ,@('Global Dialogue Forum on Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries', 'Abdallah Reyati', 'C13', 7.5, '2025-06-25 16:28:08')
    ,@('Praxislotse Wirtschaft und Menschenrechte', 'Lisa Müller', 'Undefined', 6.1, '2025-06-25 16:28:08')
    ,@('Verité Commodity Atlas', 'Dr. Maria Schmidt', 'Undefined', 8.2, '2025-06-25 16:28:08')
    ,@('Business & Human Rights Navigator', 'Prof. Dr. James Allen', 'Undefined', 6, '2025-06-25 16:28:08')
    ,@('Losing Ground, The Human Rights Impacts of Oil Palm Plantation Expansion in Indonesia', 'Abdallah Reyati', 'B', 8.1, '2025-06-25 16:28:08')
    ,@('When We Lost the Forest, We Lost Everything: Oil Palm Plantations and Rights Violations in Indonesia', 'Lisa Müller', 'B', 5, '2025-06-25 16:28:08')
)

$startRow = 36
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
